{"js": "// Fix em alguns dos artefatos\n// 1) \"Sensor Vibra\u00e7\u00e3o do Motor\" -> \"Sensor de Vibra\u00e7\u00e3o do Motor\"\n// 2) \"Comunica\u00e7\u00e3o por Bluetooth e/ou Gateway\" -> \"Comunica\u00e7\u00e3o por Bluetooth, Gateway, WiFi, etc\"\n\nconst body = context.document.body;\n\n// --- Edit 1: insert the missing \"de \" right before \"Vibra\u00e7\u00e3o\" -----------\nconst vibRanges = body.search(\"Vibra\u00e7\u00e3o\", { matchCase: true });\nvibRanges.load(\"text\");\nawait context.sync();\n\nif (vibRanges.items.length > 0) {\n  vibRanges.items[0].insertText(\"de \", \"Before\");\n  await context.sync();\n}\n\n// --- Edit 2: expand the Bluetooth bullet with Gateway/WiFi/etc ----------\nconst commRanges = body.search(\"Comunica\u00e7\u00e3o por Bluetooth e/ou Gateway\", { matchCase: true });\ncommRanges.load(\"text\");\nawait context.sync();\n\nif (commRanges.items.length > 0) {\n  commRanges.items[0].insertText(\n    \"Comunica\u00e7\u00e3o por Bluetooth, Gateway, WiFi, etc\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n", "ps1": "# Fix em alguns dos artefatos\n# 1) \"Sensor Vibra\u00e7\u00e3o do Motor\" -> \"Sensor de Vibra\u00e7\u00e3o do Motor\"\n# 2) \"Comunica\u00e7\u00e3o por Bluetooth e/ou Gateway\" -> \"Comunica\u00e7\u00e3o por Bluetooth, Gateway, WiFi, etc\"\n\n$d = $word.ActiveDocument\n\n# --- Edit 1: insert the missing \"de \" right before \"Vibra\u00e7\u00e3o\" -----------\n$rng1 = $d.Content\n$found1 = $rng1.Find.Execute(\"Vibra\u00e7\u00e3o\")\nif ($found1) {\n    $rng1.InsertBefore(\"de \")\n}\n\n# --- Edit 2: expand the Bluetooth bullet with Gateway/WiFi/etc ----------\n$rng2 = $d.Content\n$found2 = $rng2.Find.Execute(\"Comunica\u00e7\u00e3o por Bluetooth e/ou Gateway\")\nif ($found2) {\n    $rng2.Text = \"Comunica\u00e7\u00e3o por Bluetooth, Gateway, WiFi, etc\"\n}\n"}
